$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "      ""Bahnar"":[""two ara""]`n"
$ws.Range("A2").Value = "      ""Bahnar"":[""one ara""]`n"

# Avoid Excel's automatic row-height autofit (triggered by the embedded
# newline) from stamping explicit ht/customHeight attributes on the rows -
# match the target which keeps default row heights.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
